# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the aggregated "全部类型" sheet, mirroring the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F6").Value = 2811
$wsExpo.Range("F10").Value = 1585
$wsExpo.Range("F25").Value = 28
$wsExpo.Range("F27").Value = 1795

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 2811
$wsAll.Range("F11").Value = 1585
$wsAll.Range("F26").Value = 28
$wsAll.Range("F28").Value = 1795

$wb.Save()
